# Updates crypto price / volume(1h) values to match the latest scrape.
# Generated from the commit diff: 76 cell updates across rows 2-51 (columns D and E).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.849.84'
$ws.Range('E2').Value = '  +0.34%  '
$ws.Range('D3').Value = '1.625.48'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.45'
$ws.Range('E5').Value = '  -0.39%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  -0.43%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.65'
$ws.Range('E10').Value = '  +0.72%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0785'
$ws.Range('E11').Value = '  -1.02%  '
$ws.Range('D12').Value = '1.850.22'
$ws.Range('E12').Value = '  -0.39%  '
$ws.Range('E13').Value = '  -0.48%  '
$ws.Range('D14').Value = '1.617.25'
$ws.Range('E14').Value = '  -0.95%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.543'
$ws.Range('E15').Value = '  -2.67%  '
$ws.Range('E16').Value = '  -0.91%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.53'
$ws.Range('E17').Value = '  -0.87%  '
$ws.Range('D18').Value = '25.836.97'
$ws.Range('E18').Value = '  +0.26%  '
$ws.Range('E19').Value = '  -0.13%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '192.58'
$ws.Range('E20').Value = '  +0.10%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.35'
$ws.Range('E21').Value = '  -2.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.94'
$ws.Range('E22').Value = '  +0.13%  '
$ws.Range('E23').Value = '  -0.69%  '
$ws.Range('E24').Value = '  -1.50%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.999'
$ws.Range('E25').Value = '  -0.22%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '141.70'
$ws.Range('E26').Value = '  -0.78%  '
$ws.Range('E27').Value = '  +0.84%  '
$ws.Range('E28').Value = '  -0.39%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.42'
$ws.Range('E29').Value = '  -0.37%  '
$ws.Range('E30').Value = '  -0.27%  '
$ws.Range('E31').Value = '  +0.97%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.32'
$ws.Range('E32').Value = '  -0.74%  '
$ws.Range('E33').Value = '  -0.46%  '
$ws.Range('E34').Value = '  +0.42%  '
$ws.Range('E35').Value = '  +0.95%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.901'
$ws.Range('E36').Value = '  -0.13%  '
$ws.Range('D37').Value = '1.127.29'
$ws.Range('E37').Value = '  -0.61%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.546'
$ws.Range('E38').Value = '  +0.41%  '
$ws.Range('E39').Value = '  -1.91%  '
$ws.Range('E40').Value = '  +0.70%  '
$ws.Range('E41').Value = '  -0.22%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '99.44'
$ws.Range('E42').Value = '  -1.24%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.44'
$ws.Range('E43').Value = '  -1.85%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.794'
$ws.Range('E44').Value = '  -0.22%  '
$ws.Range('D45').Value = '1.761.29'
$ws.Range('E45').Value = '  -0.32%  '
$ws.Range('E46').Value = '  -1.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '56.22'
$ws.Range('E47').Value = '  +1.49%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0526'
$ws.Range('E48').Value = '  +3.89%  '
$ws.Range('E49').Value = '  +2.44%  '
$ws.Range('E50').Value = '  -0.81%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.57'
$ws.Range('E51').Value = '  +1.28%  '
